$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 546.7143
$ws.Cells.Item(18, 9).Value = 546.7143
$ws.Cells.Item(18, 11).Value = 546.7143
$ws.Cells.Item(18, 13).Value = -262.7143

$ws.Cells.Item(40, 8).Value = 2259
$ws.Cells.Item(40, 9).Value = 2443
$ws.Cells.Item(40, 10).Value = 1927.8
$ws.Cells.Item(40, 11).Value = 2443
$ws.Cells.Item(40, 12).Value = 1927.8
$ws.Cells.Item(40, 13).Value = -2268
$ws.Cells.Item(40, 14).Value = -2277.8

$ws.Cells.Item(64, 8).Value = 9468.166999999999
$ws.Cells.Item(64, 9).Value = 10521.6
$ws.Cells.Item(64, 10).Value = 4201
$ws.Cells.Item(64, 11).Value = 10521.6
$ws.Cells.Item(64, 12).Value = 4201
$ws.Cells.Item(64, 13).Value = -10273.6
$ws.Cells.Item(64, 14).Value = -4697

$ws.Cells.Item(67, 8).Value = 9468.166999999999
$ws.Cells.Item(67, 9).Value = 10521.6
$ws.Cells.Item(67, 10).Value = 4201
$ws.Cells.Item(67, 11).Value = 10521.6
$ws.Cells.Item(67, 12).Value = 4201
$ws.Cells.Item(67, 13).Value = -9663.6
$ws.Cells.Item(67, 14).Value = -5917

$ws.Cells.Item(70, 8).Value = 5137.7827
$ws.Cells.Item(70, 9).Value = 5150
$ws.Cells.Item(70, 10).Value = 5134.3887
$ws.Cells.Item(70, 11).Value = 15450
$ws.Cells.Item(70, 12).Value = 15403.1661
$ws.Cells.Item(70, 13).Value = -15180
$ws.Cells.Item(70, 14).Value = -15943.1661

$ws.Cells.Item(73, 8).Value = 5137.7827
$ws.Cells.Item(73, 9).Value = 5150
$ws.Cells.Item(73, 10).Value = 5134.3887
$ws.Cells.Item(73, 11).Value = 15450
$ws.Cells.Item(73, 12).Value = 15403.1661
$ws.Cells.Item(73, 13).Value = -14514
$ws.Cells.Item(73, 14).Value = -17275.1661

$ws.Cells.Item(76, 8).Value = 6254911.5
$ws.Cells.Item(76, 9).Value = 8338111.5
$ws.Cells.Item(76, 10).Value = 5312
$ws.Cells.Item(76, 11).Value = 8338111.5
$ws.Cells.Item(76, 12).Value = 5312
$ws.Cells.Item(76, 13).Value = -8337796.5
$ws.Cells.Item(76, 14).Value = -5942

$ws.Cells.Item(79, 8).Value = 6254911.5
$ws.Cells.Item(79, 9).Value = 8338111.5
$ws.Cells.Item(79, 10).Value = 5312
$ws.Cells.Item(79, 11).Value = 8338111.5
$ws.Cells.Item(79, 12).Value = 5312
$ws.Cells.Item(79, 13).Value = -8337019.5
$ws.Cells.Item(79, 14).Value = -7496

$ws.Cells.Item(88, 8).Value = 6813.3335
$ws.Cells.Item(88, 9).Value = 4783.3335
$ws.Cells.Item(88, 10).Value = 8166.6665
$ws.Cells.Item(88, 11).Value = 4783.3335
$ws.Cells.Item(88, 12).Value = 8166.6665
$ws.Cells.Item(88, 13).Value = -4377.3335
$ws.Cells.Item(88, 14).Value = -8978.666499999999

$ws.Cells.Item(91, 8).Value = 6813.3335
$ws.Cells.Item(91, 9).Value = 4783.3335
$ws.Cells.Item(91, 10).Value = 8166.6665
$ws.Cells.Item(91, 11).Value = 4783.3335
$ws.Cells.Item(91, 12).Value = 8166.6665
$ws.Cells.Item(91, 13).Value = -3379.3335
$ws.Cells.Item(91, 14).Value = -10974.6665

$ws.Cells.Item(98, 8).Value = 5692.5
$ws.Cells.Item(98, 9).Value = 6301.5713
$ws.Cells.Item(98, 11).Value = 6301.5713
$ws.Cells.Item(98, 13).Value = -4803.5713

$ws.Cells.Item(100, 8).Value = 8218.75
$ws.Cells.Item(100, 9).Value = 1979.1666
$ws.Cells.Item(100, 10).Value = 9466.666999999999
$ws.Cells.Item(100, 11).Value = 1979.1666
$ws.Cells.Item(100, 12).Value = 9466.666999999999
$ws.Cells.Item(100, 13).Value = -1438.1666
$ws.Cells.Item(100, 14).Value = -10548.667

$ws.Cells.Item(122, 8).Value = 5692.5
$ws.Cells.Item(122, 9).Value = 6301.5713
$ws.Cells.Item(122, 11).Value = 18904.7139
$ws.Cells.Item(122, 13).Value = -16454.7139

$ws.Cells.Item(137, 8).Value = 27789380
$ws.Cells.Item(137, 9).Value = 83334340
$ws.Cells.Item(137, 11).Value = 250003020
$ws.Cells.Item(137, 13).Value = -250000470

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 11823.889
$ws.Cells.Item(45, 10).Value = 2166.3333
$ws.Cells.Item(45, 12).Value = 2166.3333
$ws.Cells.Item(45, 14).Value = -2920.3333

$ws.Cells.Item(74, 8).Value = 328544.75
$ws.Cells.Item(74, 9).Value = 1251500
$ws.Cells.Item(74, 11).Value = 1251500
$ws.Cells.Item(74, 13).Value = -1250626

$ws.Cells.Item(77, 8).Value = 328544.75
$ws.Cells.Item(77, 9).Value = 1251500
$ws.Cells.Item(77, 11).Value = 6257500
$ws.Cells.Item(77, 13).Value = -6253132

$ws.Cells.Item(97, 8).Value = 2058468.9
$ws.Cells.Item(97, 9).Value = 2849895.5
$ws.Cells.Item(97, 11).Value = 2849895.5
$ws.Cells.Item(97, 13).Value = -2849399.5

$ws.Cells.Item(132, 8).Value = 3617.1667
$ws.Cells.Item(132, 9).Value = 2038.1875
$ws.Cells.Item(132, 11).Value = 6114.5625
$ws.Cells.Item(132, 13).Value = -3584.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1181.8
$ws.Cells.Item(22, 9).Value = 1297.5555
$ws.Cells.Item(22, 10).Value = 140
$ws.Cells.Item(22, 11).Value = 1297.5555
$ws.Cells.Item(22, 12).Value = 140
$ws.Cells.Item(22, 13).Value = -1124.5555
$ws.Cells.Item(22, 14).Value = -486

$ws.Cells.Item(94, 8).Value = 1478.7368
$ws.Cells.Item(94, 9).Value = 1501.32
$ws.Cells.Item(94, 11).Value = 1501.32
$ws.Cells.Item(94, 13).Value = -1050.32

$ws.Cells.Item(134, 8).Value = 14625.3125
$ws.Cells.Item(134, 9).Value = 20615.143
$ws.Cells.Item(134, 11).Value = 61845.429
$ws.Cells.Item(134, 13).Value = -59310.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 50004484
$ws.Cells.Item(31, 9).Value = 125001064
$ws.Cells.Item(31, 11).Value = 125001064
$ws.Cells.Item(31, 13).Value = -125000769

$ws.Cells.Item(34, 8).Value = 50004484
$ws.Cells.Item(34, 9).Value = 125001064
$ws.Cells.Item(34, 11).Value = 125001064
$ws.Cells.Item(34, 13).Value = -125000862

$ws.Cells.Item(62, 8).Value = 8272.4
$ws.Cells.Item(62, 10).Value = 10032.857
$ws.Cells.Item(62, 12).Value = 10032.857
$ws.Cells.Item(62, 14).Value = -11280.857

$ws.Cells.Item(65, 8).Value = 8272.4
$ws.Cells.Item(65, 10).Value = 10032.857
$ws.Cells.Item(65, 12).Value = 50164.285
$ws.Cells.Item(65, 14).Value = -56404.285

$ws.Cells.Item(122, 8).Value = 100783.2
$ws.Cells.Item(122, 9).Value = 143584.58
$ws.Cells.Item(122, 10).Value = 913.3333
$ws.Cells.Item(122, 11).Value = 430753.74
$ws.Cells.Item(122, 12).Value = 2739.9999
$ws.Cells.Item(122, 13).Value = -428303.74
$ws.Cells.Item(122, 14).Value = -7639.9999

$ws.Cells.Item(132, 8).Value = 37373.297
$ws.Cells.Item(132, 9).Value = 3516.7036
$ws.Cells.Item(132, 11).Value = 10550.1108
$ws.Cells.Item(132, 13).Value = -8020.110799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 30470814
$ws.Cells.Item(4, 9).Value = 44029280
$ws.Cells.Item(4, 11).Value = 132087840
$ws.Cells.Item(4, 13).Value = -132087728

$ws.Cells.Item(64, 8).Value = 1000
$ws.Cells.Item(64, 9).Value = 1000
$ws.Cells.Item(64, 11).Value = 3000
$ws.Cells.Item(64, 13).Value = -2730

$ws.Cells.Item(67, 8).Value = 1000
$ws.Cells.Item(67, 9).Value = 1000
$ws.Cells.Item(67, 11).Value = 3000
$ws.Cells.Item(67, 13).Value = -2064

$ws.Cells.Item(88, 8).Value = 3999.125
$ws.Cells.Item(88, 9).Value = 3999.125
$ws.Cells.Item(88, 11).Value = 11997.375
$ws.Cells.Item(88, 13).Value = -11569.375

$ws.Cells.Item(91, 8).Value = 3999.125
$ws.Cells.Item(91, 9).Value = 3999.125
$ws.Cells.Item(91, 11).Value = 11997.375
$ws.Cells.Item(91, 13).Value = -10515.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 29231.25
$ws.Cells.Item(52, 10).Value = 29231.25
$ws.Cells.Item(52, 12).Value = 29231.25
$ws.Cells.Item(52, 14).Value = -29749.25

$ws.Cells.Item(97, 8).Value = 1881.3334
$ws.Cells.Item(97, 9).Value = 1957.6
$ws.Cells.Item(97, 11).Value = 1957.6
$ws.Cells.Item(97, 13).Value = -1461.6

$ws.Cells.Item(126, 8).Value = 2906.0356
$ws.Cells.Item(126, 10).Value = 4156.625
$ws.Cells.Item(126, 12).Value = 12469.875
$ws.Cells.Item(126, 14).Value = -17409.875

$ws.Cells.Item(132, 8).Value = 6174.25
$ws.Cells.Item(132, 9).Value = 3199
$ws.Cells.Item(132, 10).Value = 12719.8
$ws.Cells.Item(132, 11).Value = 9597
$ws.Cells.Item(132, 12).Value = 38159.39999999999
$ws.Cells.Item(132, 13).Value = -7067
$ws.Cells.Item(132, 14).Value = -43219.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2633.72
$ws.Cells.Item(22, 10).Value = 3533.8333
$ws.Cells.Item(22, 12).Value = 3533.8333
$ws.Cells.Item(22, 14).Value = -4123.8333

$ws.Cells.Item(27, 8).Value = 2633.72
$ws.Cells.Item(27, 10).Value = 3533.8333
$ws.Cells.Item(27, 12).Value = 3533.8333
$ws.Cells.Item(27, 14).Value = -3747.8333

$ws.Cells.Item(45, 8).Value = 25020.5
$ws.Cells.Item(45, 9).Value = 30041
$ws.Cells.Item(45, 10).Value = 20000
$ws.Cells.Item(45, 11).Value = 30041
$ws.Cells.Item(45, 12).Value = 20000
$ws.Cells.Item(45, 13).Value = -29634
$ws.Cells.Item(45, 14).Value = -20814

$ws.Cells.Item(46, 8).Value = 9261.166999999999
$ws.Cells.Item(46, 9).Value = 2999.5
$ws.Cells.Item(46, 10).Value = 10043.875
$ws.Cells.Item(46, 11).Value = 2999.5
$ws.Cells.Item(46, 12).Value = 10043.875
$ws.Cells.Item(46, 13).Value = -2811.5
$ws.Cells.Item(46, 14).Value = -10419.875

$ws.Cells.Item(48, 8).Value = 31680
$ws.Cells.Item(48, 10).Value = 35000
$ws.Cells.Item(48, 12).Value = 35000
$ws.Cells.Item(48, 14).Value = -36322

$ws.Cells.Item(55, 8).Value = 587.75
$ws.Cells.Item(55, 9).Value = 580.3
$ws.Cells.Item(55, 11).Value = 580.3
$ws.Cells.Item(55, 13).Value = -407.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 10
$ws.Cells.Item(9, 9).Value = 10
$ws.Cells.Item(9, 11).Value = 10
$ws.Cells.Item(9, 13).Value = 130

$ws.Cells.Item(81, 8).Value = 11502.132
$ws.Cells.Item(81, 9).Value = 1472.0667
$ws.Cells.Item(81, 11).Value = 2944.1334
$ws.Cells.Item(81, 13).Value = -1883.1334

$ws.Cells.Item(84, 8).Value = 11502.132
$ws.Cells.Item(84, 9).Value = 1472.0667
$ws.Cells.Item(84, 11).Value = 14720.667
$ws.Cells.Item(84, 13).Value = -9416.667000000001

$ws.Cells.Item(100, 8).Value = 5256.1816
$ws.Cells.Item(100, 9).Value = 5256.1816
$ws.Cells.Item(100, 11).Value = 10512.3632
$ws.Cells.Item(100, 13).Value = -9971.3632

$ws.Cells.Item(107, 8).Value = 3336.6428
$ws.Cells.Item(107, 9).Value = 2880.875
$ws.Cells.Item(107, 11).Value = 8642.625
$ws.Cells.Item(107, 13).Value = -6722.625

$ws.Cells.Item(126, 8).Value = 63640.94
$ws.Cells.Item(126, 9).Value = 70860.13
$ws.Cells.Item(126, 11).Value = 212580.39
$ws.Cells.Item(126, 13).Value = -210110.39
